# Auto commit update: refresh Metrics figures and reposition active cells.
$wb = $excel.ActiveWorkbook
$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsToday = $wb.Worksheets.Item("today")

# Updated source figures on the Metrics sheet (column B, rows 2-13).
$wsMetrics.Range("B2").Value = 372820.87000000005
$wsMetrics.Range("B3").Value = 319183.36999999994
$wsMetrics.Range("B4").Value = 114081.56999999999
$wsMetrics.Range("B5").Value = 15125
$wsMetrics.Range("B6").Value = 5575527.9800000004
$wsMetrics.Range("B7").Value = 4719536.33
$wsMetrics.Range("B8").Value = 1646038.4500000002
$wsMetrics.Range("B9").Value = 217832
$wsMetrics.Range("B10").Value = 34040908.969999999
$wsMetrics.Range("B11").Value = 31994811.490000002
$wsMetrics.Range("B12").Value = 11927760.489999995
$wsMetrics.Range("B13").Value = 1315462

# All dependent formulas (today!B11:B22, E11:E22, F11:F22 and the
# TODAY()-1 cell in today!A1) recalculate automatically from these.

# Restore the recorded cell selections for each sheet. "today" is
# activated last so it remains the workbook's selected/visible tab.
$wsMetrics.Activate()
$wsMetrics.Range("D23").Select()

$wsToday.Activate()
$wsToday.Range("D8").Select()
